$p = $ppt.ActivePresentation
$s = $p.Slides.Item(32)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(4, 1)

# Select the "use Java’s " span within the paragraph and retype it as
# "use Kotlin’s " so that PowerPoint splits the run exactly like a manual edit would.
$target = $para.Characters(31, 11)
$target.Text = "use Kotlin’s "
